$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row: "_old" columns become "_FV2404", "_new" columns become "_FV2410".
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2. Turn the data range into an Excel Table ("Table1") with the renamed headers.
$rng = $ws.Range("A1:U71")
$lo = $ws.ListObjects.Add(1, $rng, $true, $null, $null)
$lo.Name = "Table1"

# 3. Freeze the header row (freeze pane below row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
